$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "unnamed: 1_level_1" header to "total"
$ws.Range("B2").Value = "total"

# Remove the two label-only section-header rows that had no data next to
# them ("situação do domicílio" and "grandes regiões e unidades da
# federação"), which fixes the one-row data/label misalignment below each
# of them. Row 5 is the first one; once it's removed the second one (which
# was row 8) becomes row 7.
$ws.Rows(5).Delete()
$ws.Rows(7).Delete()
